$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regex strings (password/email regexes now escape "/" as "\/";
# also a new copy of the email regex is placed at row 14, the bcrypt regex
# at row 8, and the password regex at row 12 - per the authoritative diff).
$bcryptRegex = '^[$]2[abxy]?[$](?:0[4-9]|[12][0-9]|3[01])[$][.\/0-9a-zA-Z]{53}$'
$passwordRegex = '(?=(.*[0-9]))(?=.*[\!@#$%^&*()\\[\]{}\-_+=~`|:;"''<>,.\/?])(?=.*[a-z])(?=(.*[A-Z]))(?=(.*)).{8,}'
$emailRegex = '(?:[a-z0-9!#$%&''*+\/=?^_`{|}~-]+(?:\.[a-z0-9!#$%&''*+\/=?^_`{|}~-]+)*|"(?:[\x01-\x08\x0b\x0c\x0e-\x1f\x21\x23-\x5b\x5d-\x7f]|\\[\x01-\x09\x0b\x0c\x0e-\x7f])*")@(?:(?:[a-z0-9](?:[a-z0-9-]*[a-z0-9])?\.)+[a-z0-9](?:[a-z0-9-]*[a-z0-9])?|\[(?:(?:25[0-5]|2[0-4][0-9]|[01]?[0-9][0-9]?)\.){3}(?:25[0-5]|2[0-4][0-9]|[01]?[0-9][0-9]?|[a-z0-9-]*[a-z0-9]:(?:[\x01-\x08\x0b\x0c\x0e-\x1f\x21-\x5a\x53-\x7f]|\\[\x01-\x09\x0b\x0c\x0e-\x7f])+)\])'
$phoneNumberName = 'phoneNumber'
$phoneNumberRegex = '^\+[1-9]\d{1,14}$'
$bcryptHashName = 'bcrypt hash'
$unsignedIntRegex = '^\d+$'
$signedFloatRegex = '[-+]?([0-9]*[.])?[0-9]+([eE][-+]?\d+)?'
$trueFalseRegex = '(?:true|false)'
$yesNoRegex = '(?:yes|no)'
$onOffRegex = '(?:on|off)'

# Row 8 (email) now carries the bcrypt-hash format regex.
$ws.Range("C8").Value = $bcryptRegex

# Row 13 becomes phoneNumber / phoneNumber-regex.
$ws.Range("B13").Value = $phoneNumberName
$ws.Range("C13").Value = $phoneNumberRegex

# Row 14 becomes "bcrypt hash" / email-regex. (Write the email regex before
# the password regex below so the shared-string table regenerates in the
# same order as the authored workbook.)
$ws.Range("B14").Value = $bcryptHashName
$ws.Range("C14").Value = $emailRegex

# Row 12 (password) now carries the updated password regex (escaped "/").
$ws.Range("C12").Value = $passwordRegex

# Rows 16/18/20/22 (the "*unsigned" rows) switch from "^\d+$" (shared with
# row 16's old sibling) to the de-duplicated "^\d+$" string.
$ws.Range("C16").Value = $unsignedIntRegex
$ws.Range("C18").Value = $unsignedIntRegex
$ws.Range("C20").Value = $unsignedIntRegex
$ws.Range("C22").Value = $unsignedIntRegex

# Rows 23/24 (float/double).
$ws.Range("C23").Value = $signedFloatRegex
$ws.Range("C24").Value = $signedFloatRegex

# Rows 25/26/27 (true-false / yes-no / on-off).
$ws.Range("C25").Value = $trueFalseRegex
$ws.Range("C26").Value = $yesNoRegex
$ws.Range("C27").Value = $onOffRegex

# Update the last-used selection to match the authored file.
$ws.Range("I21").Select() | Out-Null
